$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.196.03'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.27%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.429.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.88%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '407.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.96%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.99%  '

$ws.Range('E7').Value = '  -1.41%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.694'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.77%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.134'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.05'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.60%  '

$ws.Range('E12').Value = '  +0.08%  '

$ws.Range('E13').Value = '  +1.62%  '

$ws.Range('E14').Value = '  -0.72%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.430.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '11.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.26%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.152.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.28%  '

$ws.Range('E18').Value = '  -0.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000149'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.97%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '84.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.84%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '311.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.00%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.81%  '

$ws.Range('E24').Value = '  +1.14%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.99%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.75'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.83%  '

$ws.Range('E30').Value = '  +0.16%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '43.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.16%  '

$ws.Range('E32').Value = '  -0.46%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.37'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.36%  '

$ws.Range('E34').Value = '  +0.10%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0486'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.64'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('E38').Value = '  +1.46%  '

$ws.Range('E39').Value = '  -2.37%  '

$ws.Range('E40').Value = '  +12.50%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.42%  '

$ws.Range('E42').Value = '  +0.29%  '

$ws.Range('E43').Value = '  -2.26%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.60%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.35%  '

$ws.Range('E46').Value = '  +0.18%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.89%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.105.62'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.91%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.58%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +19.64%  '
